$d = $word.ActiveDocument

# 1. Fix capitalization: "Github Username" -> "GitHub Username" in the signature table.
$d.Content.Find.Execute("Github Username", $true, $false, $false, $false, $false, $true, 1, $false, "GitHub Username", 2)

# 2. Remove the trailing "------" divider, the blank paragraph after it, and the
#    "If you have any attachments..." paragraph that followed the signature table.
$count = $d.Paragraphs.Count
$target = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^------") {
        $target = $i
        break
    }
}
if ($target -gt 0) {
    $startPara = $d.Paragraphs.Item($target)
    $endPara = $d.Paragraphs.Item($target + 2)
    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRange.Delete()
}

# 3. Footer updates: bump the CLA version number and refresh the "Last Updated" date.
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)

$r1 = $ftr.Range
$r1.Find.Execute("Version 1.1", $true, $false, $false, $false, $false, $true, 1, $false, "Version 1.2", 2)

$r2 = $ftr.Range
$r2.Find.Execute("Last Updated: March 3rd, 2022", $true, $false, $false, $false, $false, $true, 1, $false, "Last Updated: April 26th, 2023", 2)
